$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.400.82"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "1.916.20"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4691"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2849"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06811"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "107.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("D12").Value = "1.898.11"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07629"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6538"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "288.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").Value = "30.408.10"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007620"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "2.148.98"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.221"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.197"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.255"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.035"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1068"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.139"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.934"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05022"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7371"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.145"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +6.78%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.045"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8730"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.847"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4202"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "52.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +25.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.161"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.204"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1205"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3872"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.57%  "
